$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "85.963.29"
$c.ClearFormats()
$ws.Range("E2").Value = "  +4.66%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.264.08"
$c.ClearFormats()
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  +0.50%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "210.26"
$c.ClearFormats()
$ws.Range("E5").Value = "  -2.73%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "618.65"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.40%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.368"
$c.ClearFormats()
$ws.Range("E7").Value = "  +27.64%  "
$ws.Range("E8").Value = "  +0.25%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.624"
$c.ClearFormats()
$ws.Range("E9").Value = "  +6.35%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.262.81"
$c.ClearFormats()
$ws.Range("E10").Value = "  +2.13%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.577"
$c.ClearFormats()
$ws.Range("E11").Value = "  -2.19%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0000257"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.99%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.173"
$c.ClearFormats()
$ws.Range("E13").Value = "  +5.05%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.899.66"
$c.ClearFormats()
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "34.12"
$c.ClearFormats()
$ws.Range("E15").Value = "  +7.32%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.18"
$c.ClearFormats()
$ws.Range("E16").Value = "  -2.75%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "86.021.30"
$c.ClearFormats()
$ws.Range("E17").Value = "  +5.02%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.288.39"
$c.ClearFormats()
$ws.Range("E18").Value = "  +3.02%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.05"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.11%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "9.05"
$c.ClearFormats()
$ws.Range("E20").Value = "  +0.82%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "431.05"
$c.ClearFormats()
$ws.Range("E21").Value = "  -1.15%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.ClearFormats()
$ws.Range("E22").Value = "  -9.96%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.ClearFormats()
$ws.Range("E23").Value = "  +3.44%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "7.22"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.67%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "5.14"
$c.ClearFormats()
$ws.Range("E25").Value = "  -2.74%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "11.95"
$c.ClearFormats()
$ws.Range("E26").Value = "  +8.96%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "3.473.33"
$c.ClearFormats()
$ws.Range("E27").Value = "  +3.70%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "76.48"
$c.ClearFormats()
$ws.Range("E28").Value = "  -0.40%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.0000127"
$c.ClearFormats()
$ws.Range("E29").Value = "  +2.69%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$ws.Range("E30").Value = "  +0.00%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.169"
$c.ClearFormats()
$ws.Range("E31").Value = "  +22.00%  "
$ws.Range("E32").Value = "  +0.19%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "8.77"
$c.ClearFormats()
$ws.Range("E33").Value = "  -3.12%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "539.93"
$c.ClearFormats()
$ws.Range("E34").Value = "  -8.50%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.43"
$c.ClearFormats()
$ws.Range("E35").Value = "  -5.56%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.98"
$c.ClearFormats()
$ws.Range("E36").Value = "  -1.32%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.63"
$c.ClearFormats()
$ws.Range("E37").Value = "  +8.03%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.136"
$c.ClearFormats()
$ws.Range("E38").Value = "  -12.47%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "22.49"
$c.ClearFormats()
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("E40").Value = "  +0.36%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "21.82"
$c.ClearFormats()
$ws.Range("E41").Value = "  +4.95%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.402"
$c.ClearFormats()
$ws.Range("E42").Value = "  -1.63%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.98"
$c.ClearFormats()
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.ClearFormats()
$ws.Range("E45").Value = "  -5.57%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "153.83"
$c.ClearFormats()
$ws.Range("E46").Value = "  -4.28%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "177.75"
$c.ClearFormats()
$ws.Range("E47").Value = "  -5.61%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "44.49"
$c.ClearFormats()
$ws.Range("E48").Value = "  -0.22%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.32"
$c.ClearFormats()
$ws.Range("E49").Value = "  -0.72%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.23"
$c.ClearFormats()
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.624"
$c.ClearFormats()
$ws.Range("E51").Value = "  -0.86%  "
